$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 642: 2020-10-02 -> 18591.26
$ws.Range("A641:B641").Copy()
$ws.Range("A642:B642").PasteSpecial(-4122)
$ws.Range("A642").Value = 44106
$ws.Range("B642").Value = 18591.259999999998

# Row 643: 2020-10-03 (Saturday) -> "--"
$ws.Range("A628:B628").Copy()
$ws.Range("A643:B643").PasteSpecial(-4122)
$ws.Range("A643").Value = 44107
$ws.Range("B643").Value = "--"

# Row 644: 2020-10-04 (Sunday) -> "--"
$ws.Range("A628:B628").Copy()
$ws.Range("A644:B644").PasteSpecial(-4122)
$ws.Range("A644").Value = 44108
$ws.Range("B644").Value = "--"

# Row 645: 2020-10-05 -> 18633.89
$ws.Range("A641:B641").Copy()
$ws.Range("A645:B645").PasteSpecial(-4122)
$ws.Range("A645").Value = 44109
$ws.Range("B645").Value = 18633.89

# Row 646: 2020-10-06 -> 18398.08
$ws.Range("A641:B641").Copy()
$ws.Range("A646:B646").PasteSpecial(-4122)
$ws.Range("A646").Value = 44110
$ws.Range("B646").Value = 18398.080000000002

# Row 647: 2020-10-07 -> 18235.02
$ws.Range("A641:B641").Copy()
$ws.Range("A647:B647").PasteSpecial(-4122)
$ws.Range("A647").Value = 44111
$ws.Range("B647").Value = 18235.02

# Row 648: 2020-10-08 -> 18497.59
$ws.Range("A641:B641").Copy()
$ws.Range("A648:B648").PasteSpecial(-4122)
$ws.Range("A648").Value = 44112
$ws.Range("B648").Value = 18497.59

# Row 649: 2020-10-09 -> 18624.21
$ws.Range("A641:B641").Copy()
$ws.Range("A649:B649").PasteSpecial(-4122)
$ws.Range("A649").Value = 44113
$ws.Range("B649").Value = 18624.21

$excel.CutCopyMode = 0

# Update the defined name range to extend through the new last row
$wb.Names.Item("IGPA").RefersTo = "=IGPA!`$A`$1:`$B`$649"

# Update view: scroll window down and move the active selection
$excel.ActiveWindow.ScrollRow = 642
$ws.Range("B652").Select() | Out-Null
